# New row-by-row contents for the Naranja (orange) price table, rows 792-836.
# This reflects a weekly data refresh: 2 new rows (Fukumoto, 2022-07-11) were
# inserted at the top of the Femacal de La Calera block, pushing the existing
# rows down by two positions, with 2 rows (Navel Late, 2021-09-15) now appended
# at the bottom of the table (rows 835-836).
$rows = @(
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44753, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Primera', 240, 3000, 3500, 3223, '$/malla 13 kilos', 'Provincia de Quillota', 248, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44753, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Segunda', 239, 2400, 2700, 2528, '$/malla 13 kilos', 'Provincia de Quillota', 194, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44489, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 165, 4500, 5000, 4758, '$/malla 13 kilos', 'Provincia de Quillota', 366, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44489, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 135, 3500, 4000, 3752, '$/malla 13 kilos', 'Provincia de Quillota', 289, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44489, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Primera', 115, 4500, 5000, 4783, '$/malla 13 kilos', 'Provincia de Quillota', 368, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44489, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Segunda', 110, 3500, 4000, 3773, '$/malla 13 kilos', 'Provincia de Quillota', 290, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44659, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Primera', 50, 8000, 8000, 8000, '$/malla 13 kilos', 'Provincia de Quillota', 615, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44659, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Segunda', 50, 7000, 7000, 7000, '$/malla 13 kilos', 'Provincia de Quillota', 538, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44505, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 172, 5000, 6000, 5494, '$/malla 13 kilos', 'Provincia de Quillota', 423, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44505, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 170, 4000, 4500, 4265, '$/malla 13 kilos', 'Provincia de Quillota', 328, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44340, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Primera', 155, 7500, 8000, 7742, '$/malla 13 kilos', 'Provincia de Quillota', 596, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44340, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Segunda', 140, 6000, 6500, 6250, '$/malla 13 kilos', 'Provincia de Quillota', 481, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44326, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Primera', 138, 10000, 11000, 10507, '$/malla 13 kilos', 'Provincia de Quillota', 808, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44326, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Segunda', 137, 8000, 9000, 8489, '$/malla 13 kilos', 'Provincia de Quillota', 653, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44714, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Primera', 143, 6500, 7000, 6797, '$/malla 13 kilos', 'Provincia de Quillota', 523, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44714, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Segunda', 138, 5500, 6000, 5783, '$/malla 13 kilos', 'Provincia de Quillota', 445, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44343, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Primera', 155, 7000, 8000, 7548, '$/malla 13 kilos', 'Provincia de Quillota', 581, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44343, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Fukumoto', 'Segunda', 153, 6000, 6500, 6222, '$/malla 13 kilos', 'Provincia de Quillota', 479, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44426, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Primera', 226, 3500, 4000, 3779, '$/malla 13 kilos', 'Provincia de Quillota', 291, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44426, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Segunda', 70, 3000, 3000, 3000, '$/malla 13 kilos', 'Provincia de Quillota', 231, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44426, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'New Hall', 'Primera', 148, 3500, 4000, 3797, '$/malla 13 kilos', 'Provincia de Quillota', 292, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44426, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'New Hall', 'Segunda', 70, 3000, 3000, 3000, '$/malla 13 kilos', 'Provincia de Quillota', 231, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44259, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Primera', 56, 13000, 13000, 13000, '$/malla 13 kilos', 'Provincia de Quillota', 1000, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44259, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Segunda', 60, 12000, 12000, 12000, '$/malla 13 kilos', 'Provincia de Quillota', 923, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44259, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Tercera', 48, 10000, 10000, 10000, '$/malla 13 kilos', 'Provincia de Quillota', 769, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44238, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Primera', 180, 13000, 13000, 13000, '$/malla 13 kilos', 'Provincia de Quillota', 1000, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44376, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 125, 4500, 5000, 4760, '$/malla 13 kilos', 'Provincia de Quillota', 366, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44376, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 60, 4000, 4000, 4000, '$/malla 13 kilos', 'Provincia de Quillota', 308, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44376, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Primera', 193, 4500, 5000, 4762, '$/malla 13 kilos', 'Provincia de Quillota', 366, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44376, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'New Hall', 'Primera', 182, 4500, 5000, 4766, '$/malla 13 kilos', 'Provincia de Quillota', 367, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44376, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'New Hall', 'Segunda', 80, 4000, 4000, 4000, '$/malla 13 kilos', 'Provincia de Quillota', 308, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44622, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Primera', 123, 7000, 7500, 7272, '$/malla 13 kilos', 'Provincia de Quillota', 559, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44622, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Segunda', 68, 6000, 6000, 6000, '$/malla 13 kilos', 'Provincia de Quillota', 462, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44406, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 150, 3500, 4000, 3733, '$/malla 13 kilos', 'Provincia de Quillota', 287, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44406, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 157, 2500, 3000, 2777, '$/malla 13 kilos', 'Provincia de Quillota', 214, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44406, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Primera', 185, 3500, 4000, 3735, '$/malla 13 kilos', 'Provincia de Quillota', 287, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44406, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Segunda', 172, 2500, 3000, 2747, '$/malla 13 kilos', 'Provincia de Quillota', 211, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44406, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'New Hall', 'Primera', 170, 3500, 4000, 3765, '$/malla 13 kilos', 'Provincia de Quillota', 290, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44627, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Valencia', 'Primera', 85, 7500, 7500, 7500, '$/malla 13 kilos', 'Provincia de Quillota', 577, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44547, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 125, 7000, 7500, 7260, '$/malla 13 kilos', 'Provincia de Quillota', 558, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44547, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 128, 5000, 6000, 5547, '$/malla 13 kilos', 'Provincia de Quillota', 427, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44454, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Primera', 242, 3500, 4000, 3783, '$/malla 13 kilos', 'Provincia de Quillota', 291, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44454, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Lane Late', 'Segunda', 70, 3000, 3000, 3000, '$/malla 13 kilos', 'Provincia de Quillota', 231, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44454, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Primera', 172, 3500, 4000, 3753, '$/malla 13 kilos', 'Provincia de Quillota', 289, 13)
    ,@(3, 'Femacal de La Calera', 'Coquimbo', 44454, 5, 'Fruta', 100102, 'Cítricos', 100102005, 'Naranja', 'Navel Late', 'Segunda', 100, 3000, 3000, 3000, '$/malla 13 kilos', 'Provincia de Quillota', 231, 13)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 792
$lastExistingRow = 834   # last row that already existed before this edit

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Rows appended past the previous last row (835, 836) start out with the
# default/general number format, so column D (the date column) needs its
# date/time format applied explicitly, matching the rest of the column.
$endRow = $startRow + $rows.Count - 1
if ($endRow -gt $lastExistingRow) {
    $newRange = $ws.Range("D" + ($lastExistingRow + 1) + ":D" + $endRow)
    $newRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Updated rows 792-836. UsedRange:" ($ws.UsedRange.Address())
